$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. B5: fill in the two evidence hashes, move the selection to A3
# ---------------------------------------------------------------------
$b5 = $wb.Worksheets.Item("B5")
$b5.Range("A2").Value = "DA89BF3CB114E5958AD762BACD4C4C6A343FEAC2C0EDAB23737610D6A9768DD0"
$b5.Range("A3").Value = "7FA45E82327B42710FF4423E4CC80DE7F60EFF7EB1E9C4967809A64631A47467"

# ---------------------------------------------------------------------
# 2. B6: fill in the two evidence hashes, move the selection to A3
# ---------------------------------------------------------------------
$b6 = $wb.Worksheets.Item("B6")
$b6.Range("A2").Value = "72620E9A66C080EE14B536E210AC3F6B71C3D5A401AFA6C595A0A4590160D065"
$b6.Range("A3").Value = "34E75DD07859B6FF24E41E380727E9668E6BA46126FA85AD58AE62219DE2AB57"

# ---------------------------------------------------------------------
# 3. Add a new sheet "B7" at the end, filled in with real evidence
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$b7 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$b7.Name = "B7"
$b7.Range("A1").Value = "TxHash"
$b7.Range("A2").Value = "C449B60765751390553D11172C966921A93CBA3A59A27B085BBA9C7DCB0A9FD3"
$b7.Range("A3").Value = "500F1C46E512219F1B9D8CF812AEE703795F41FB8144496E72C257ADA5B7B7C4"

# ---------------------------------------------------------------------
# 4. Add a new sheet "B8" at the end, still the blank template text
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$b8 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$b8.Name = "B8"
$b8.Range("A1").Value = "TxHash"
$b8.Range("A2").Value = "The first Interchain NFT-Transfer TxHash"
$b8.Range("A3").Value = "The Internal Transfer TxHash on IRISnet"

# ---------------------------------------------------------------------
# 5. Fix up the selections on each of B5/B6/B7/B8 and the active sheet
# ---------------------------------------------------------------------
$b5.Activate()
[void]$b5.Range("A3").Select()

$b6.Activate()
[void]$b6.Range("A3").Select()

$b8.Activate()
[void]$b8.Range("G16").Select()

$b7.Activate()
[void]$b7.Range("H14").Select()
